# cosas-portal EMX workbook: "fixed variable names and rebuilt emx"
#
# 1) bump the package description version/date stamp
# 2) on the "attributes" sheet, the cosasportal_bench_cnv entity gets a new
#    "Phenotype" (string) attribute in place of the old row 128, and the
#    "dataProcessed" (dateTime) attribute that used to live there is
#    reinserted right after it - pushing the trailing "id" attribute row
#    down by one.

$wb = $excel.ActiveWorkbook

# --- 1. packages sheet: update the staging-tables description string ---
$wsPackages = $wb.Worksheets.Item("packages")
$wsPackages.Range("C2").Value = "Staging tables for raw data extracts (v0.9011, 2021-08-30)"

# --- 2. attributes sheet: rework the cosasportal_bench_cnv attribute rows ---
$wsAttributes = $wb.Worksheets.Item("attributes")

# Row 128 used to be "dataProcessed" (dateTime); rename it in place to the
# new "Phenotype" (string) attribute - entity + flag columns stay the same.
$wsAttributes.Range("A128").Value = "Phenotype"
$wsAttributes.Range("C128").Value = "string"

# Insert a fresh row at 129 (this pushes the old row 129 "id" attribute
# down to row 130) and populate it with the "dataProcessed" attribute that
# used to be at row 128.
$wsAttributes.Rows.Item(129).Insert()
$wsAttributes.Range("A129").Value = "dataProcessed"
$wsAttributes.Range("B129").Value = "cosasportal_bench_cnv"
$wsAttributes.Range("C129").Value = "dateTime"
$wsAttributes.Range("G129").Value = $false
$wsAttributes.Range("I129").Value = $true
$wsAttributes.Range("J129").Value = $false

Write-Output "applied cosasportal attribute fixups"
